$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 13) after the existing rows (A1:D12 -> A1:D13).
$row = 13

# Column A holds the date as literal text (matching the existing rows,
# e.g. A2 = "08/29/2025"), not as a real date serial number. Force text
# formatting before assigning so Excel doesn't auto-convert the string
# into a date, then copy the plain/default style from an existing
# unstyled data cell so the new cell ends up without any extra explicit
# style (same as A2:A5, A7:A12).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "11/07/2025"
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item(2, 1).Style

$ws.Cells.Item($row, 2).Value = 494.851999999999
$ws.Cells.Item($row, 3).Value = 0.05001495396603439
$ws.Cells.Item($row, 4).Value = 25
